$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate() | Out-Null

# Mark the Create/Read/Update/Delete test-passed columns as TRUE for every
# data row (References added + Compared Test Data => all CRUD checks pass).
$range = $ws.Range("B2:E24")
$range.Value = $true

$ws.Range("G7").Select() | Out-Null
